$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name: -> Supervised Internship (was "Supervised Work")
$ws.Range("B4").Value = "Supervised Internship"
$ws.Range("C4").Value = "Supervised Internship"

# Créditos-trabalho: 12 -> 6 (force text, matching the source sharedString type)
$ws.Range("B6").Value = "'6"
$ws.Range("C6").Value = "'6"

# Carga horária: 375 h -> 195 h
$ws.Range("B7").Value = "195 h   (    Estágio: 195 h         )"
$ws.Range("C7").Value = "195 h   (    Estágio: 195 h         )"

# Ativação: 01/01/2012 -> 01/01/2023 (force text, not a real date cell)
$ws.Range("B8").Value = "'01/01/2023"
$ws.Range("C8").Value = "'01/01/2023"

# Objetivos: professor changed
$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# Objectives: (new paragraph text added under the label)
$objectivesText = "Offer the opportunity to carry out professional training in a company or research institution, under the supervision of a professor from the Materials Engineering Department at EEL. Complement the general curricular training and psychologically and socially adapt the student to his/her future professional activity."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# Programa resumido: new value "Semestral" (previously held a stray leftover "01/01/2012")
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Short syllabus: (new paragraph text added under the label)
$shortSyllabusText = "Participation in the selection process or indication of an institution to carry out an internship. Submission of the specific work plan. Conducting the internship and delivering the internship report."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText

# Programa: stray leftover value updated from old professor name to the new activation date
# (force text, not a real date cell)
$ws.Range("B15").Value = "'01/01/2023"
$ws.Range("C15").Value = "'01/01/2023"

# Syllabus: (new paragraph text added under the label)
$syllabusText = "Student participation in the selection process of companies, research institutions or in the academic sector. The internship will be carried out under the supervision of a professor appointed by the Physical Engineering Course Committee. The content will be established in the Work Plan between the supervisor responsible for the Internship and the supervising professor. Presentation of a final report on the activities carried out in the internship."
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText
